$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels: split the ambiguous "% YOY Change" label into
# distinct Revenue and Expenses variants.
$ws.Range("I1").Value = "% YOY Expenses Change"
$ws.Range("E1").Value = "% YOY Revenue Change"

# Update the saved selection to match the author's final cursor position.
$ws.Range("G15").Select()
